$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "61.001.47"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.379.41"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "403.83"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.08"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +15.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +6.17%  "
$ws.Range("D8").Value = "3.380.67"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.670"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +7.43%  "
$ws.Range("E11").Value = "  +15.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.07"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +8.42%  "
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "3.953.75"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.53"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.61"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "3.401.19"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.64"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +11.32%  "
$ws.Range("D19").Value = "61.057.94"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.01"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000133"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +18.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.23"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.47"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +11.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.95"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +4.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "306.41"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.22"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.64"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +15.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.46"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.60"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.44"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.115"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.92"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +11.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.64"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("E35").Value = "  +6.21%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0481"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.19"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.95"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +4.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.61"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.284"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.87"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +4.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.90"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.23"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.81"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -22.88%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.138.53"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "3.724.68"
$ws.Range("E51").Value = "  -6.49%  "
